$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-21 Thursday" "2025-08-22 Friday"

Replace-Text "679÷6=113, 1" "601÷7=85, 6"
Replace-Text "217÷2=108, 1" "703÷6=117, 1"
Replace-Text "425÷6=70, 5" "931÷2=465, 1"
Replace-Text "398÷7=56, 6" "842÷3=280, 2"
Replace-Text "591÷2=295, 1" "235÷5=47, 0"

Replace-Text "639÷7=91, 2" "163÷9=18, 1"
Replace-Text "957÷5=191, 2" "174÷2=87, 0"
Replace-Text "572÷8=71, 4" "137÷3=45, 2"
Replace-Text "484÷8=60, 4" "994÷4=248, 2"
Replace-Text "276÷8=34, 4" "141÷7=20, 1"

Replace-Text "215÷2=107, 1" "738÷8=92, 2"
Replace-Text "499÷9=55, 4" "956÷4=239, 0"
Replace-Text "702÷8=87, 6" "406÷7=58, 0"
Replace-Text "269÷7=38, 3" "970÷5=194, 0"
Replace-Text "347÷6=57, 5" "973÷9=108, 1"

Replace-Text "522÷4=130, 2" "784÷2=392, 0"
Replace-Text "433÷9=48, 1" "328÷4=82, 0"
Replace-Text "913÷6=152, 1" "417÷3=139, 0"
Replace-Text "225÷8=28, 1" "822÷5=164, 2"
Replace-Text "355÷9=39, 4" "669÷9=74, 3"

Replace-Text "301÷7=43, 0" "763÷2=381, 1"
Replace-Text "973÷3=324, 1" "699÷7=99, 6"
Replace-Text "769÷4=192, 1" "473÷6=78, 5"
Replace-Text "138÷4=34, 2" "641÷7=91, 4"
Replace-Text "787÷6=131, 1" "974÷3=324, 2"
